# RPA datasets push 2023-12-06
# Two new IPO rows (교보15호스팩 / 와이바이오로직스, both listed 2023-12-05)
# were added to the top of the data table, pushing every existing row
# down by two rows (old row 2 -> new row 4, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right after the header row, shifting the
# existing data (old rows 2-31) down to rows 4-33.
$ws.Rows("2:3").Insert()

# The insert operation can drag the header row's bold/centered format
# onto the newly created rows - strip that back to the plain/default
# style used by every other data row.
$ws.Range("A2:Q3").ClearFormats()

# Columns A, O and P hold plain "yyyy-mm-dd" text labels (not real
# dates) everywhere else in the sheet. Mark them as Text before writing
# so Excel doesn't silently convert the strings into date serials.
$dateCols = @("A2:A3", "O2:O3", "P2:P3")
foreach ($rng in $dateCols) {
    $ws.Range($rng).NumberFormat = "@"
}

# New row 2: 교보15호스팩
$ws.Range("A2").Value = "2023-12-05"
$ws.Range("B2").Value = "교보15호스팩"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 70
$ws.Range("E2").Value = "교보"
$ws.Range("F2").Value = 70
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 2000
$ws.Range("N2").Value = 100
$ws.Range("O2").Value = "2023-11-23"
$ws.Range("P2").Value = "2023-11-28"
$ws.Range("Q2").Value = 2625000

# New row 3: 와이바이오로직스
$ws.Range("A3").Value = "2023-12-05"
$ws.Range("B3").Value = "와이바이오로직스"
$ws.Range("C3").Value = "코스닥"
$ws.Range("D3").Value = 135
$ws.Range("E3").Value = "유안타"
$ws.Range("F3").Value = 135
$ws.Range("G3").Value = "-"
$ws.Range("H3").Value = "-"
$ws.Range("I3").Value = "-"
$ws.Range("J3").Value = "-"
$ws.Range("K3").Value = "대표"
$ws.Range("L3").Value = "-"
$ws.Range("M3").Value = 9000
$ws.Range("N3").Value = 100
$ws.Range("O3").Value = "2023-11-23"
$ws.Range("P3").Value = "2023-11-28"
$ws.Range("Q3").Value = 1055000

# Put those three columns back on the default/plain style (no explicit
# number format override), matching every other data row in the sheet.
foreach ($rng in $dateCols) {
    $ws.Range($rng).Style = "Normal"
}
